# Fix bug related to events with the same ID, add CM groups to course info
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("L1")

# Insert a new column before the existing "Heures TD" column (column F),
# shifting Heures TD / Groupes TD / Heures TP / Groupes TP one column to
# the right, and populate it with the new "Groupes CM" header + values.
$ws.Columns("F").Insert()

$ws.Range("F1").Value2 = "Groupes CM"

$groupesCM = @{
    2  = 1
    3  = 1
    4  = 3.5
    5  = 1
    6  = 1
    7  = 1
    8  = 3.5
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 1
}

foreach ($row in $groupesCM.Keys) {
    $ws.Cells.Item($row, 6).Value2 = $groupesCM[$row]
}

# Refresh the sort-state range so it covers the newly inserted column.
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("B2:B31"))
$sortObj.SortFields.Add($ws.Range("A2:A31"))
$sortObj.SetRange($ws.Range("A2:J31"))
$sortObj.Header = 0
$sortObj.Apply()

# Update the active selection to match the saved view state.
$ws.Select()
$ws.Range("C21").Select()
